# Update "想去人数" (interested-count) figures scraped for two events.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 652
$wsExpo.Range("F4").Value = 1494
$wsExpo.Range("F5").Value = 697

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 652
$wsAll.Range("F4").Value = 1494
$wsAll.Range("F6").Value = 697
